# Applies the change described by the diff: the worksheet's row 10
# (hasGeometry / GeomValue / Geometric object / Geometrisches Objekt)
# is removed entirely, shifting all subsequent rows up by one.
#
# Deleting the whole row via the Excel object model automatically:
#  - shifts rows 11..34 up to 10..33
#  - recomputes the sheet <dimension> (A1:S34 -> A1:S33)
#  - drops the now-unused shared strings ("hasGeometry", "GeomValue",
#    "Geometric object", "Geometrisches Objekt") when Excel re-saves
#  - updates the sheet's active cell/selection bookkeeping

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Rows.Item(10).Delete()

# Mirror Excel's own post-delete selection (whole new row 10 selected).
$ws.Range("A10:XFD10").Select()

$wb.Save()
